$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CasesTab query (B2): append a new line with an ORDER BY / LIMIT clause ---
$b2cur = $ws.Range("B2").Value2
$b2suffix = @'

 order By ss.study_subject_id ASC LIMIT 100
'@
$ws.Range("B2").Value = $b2cur + $b2suffix

# --- SamplesTab query (B3): append a new line with an ORDER BY / LIMIT clause ---
$b3cur = $ws.Range("B3").Value2
$b3suffix = @'

order By samp.sample_id ASC LIMIT 100
'@
$ws.Range("B3").Value = $b3cur + $b3suffix

# --- FilesTab query (B4): replace the trailing "order by" clause in place ---
$b4cur = $ws.Range("B4").Value2
$b4old = "    order by f.file_name"
$b4new = "     order By f.file_name ASC LIMIT 100"
$ws.Range("B4").Value = $b4cur.Substring(0, $b4cur.Length - $b4old.Length) + $b4new

# --- Row heights: text grew, so the (auto) wrapped-row heights grew too ---
$ws.Rows("2").RowHeight = 360
$ws.Rows("3").RowHeight = 360

# --- Selection moved from C4 to B4 (with the view scrolled down a bit) ---
$ws.Range("B4").Select()
